$wb = $excel.ActiveWorkbook

# The "Meta" sheet header in B1 was renamed from the long assembly-qualified
# type name label to simply "Type".
$metaSheet = $wb.Worksheets.Item("Meta")
$metaSheet.Range("B1").Value = "Type"

# Make "Meta" the active sheet/tab again (it was "Resources" before), with
# B1 selected.
$metaSheet.Activate()
$metaSheet.Range("B1").Select()
